# Applies the edits described by the target diff:
#  - Gesamtinvestitionskosten sheet: update several "netto" input cells (column B)
#    and the Steuerkorrektur rate (B20), plus one literal percentage (D10).
#    All dependent formulas (C/D/E/F columns, totals in rows 12/14) recalculate
#    automatically.
#  - Update the remembered cell-selection (active cell) on three sheets, matching
#    where the user last clicked while testing the new input validation.

$wb = $excel.ActiveWorkbook

# --- Sheet "Gesamtinvestitionskosten" -------------------------------------
$ws = $wb.Worksheets.Item("Gesamtinvestitionskosten")

$ws.Range("B2").Value  = 10
$ws.Range("B3").Value  = 20
$ws.Range("B4").Value  = 8
$ws.Range("B5").Value  = 10
$ws.Range("B6").Value  = 5
$ws.Range("B8").Value  = 643
$ws.Range("D10").Value = 0.7
$ws.Range("B20").Value = 0.5

$ws.Range("B21").Select()

# --- Sheet "Basisinformation" ----------------------------------------------
$ws1 = $wb.Worksheets.Item("Basisinformation")
$ws1.Activate()
$ws1.Range("F14").Select()

# --- Sheet "Wirtschaftlichkeitsrechnung" -----------------------------------
$ws4 = $wb.Worksheets.Item("Wirtschaftlichkeitsrechnung")
$ws4.Activate()
$ws4.Range("E20").Select()

# Leave the originally active sheet ("Gesamtinvestitionskosten") selected/active,
# matching tabSelected="1" in the workbook.
$ws.Activate()
$ws.Range("B21").Select()

$wb.Save()
